# Initial touch pad code
# Adds new translation rows (31-40) to the "Translation" sheet and
# updates a few existing cells (F11, C29, D29) on that same sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Update existing cells ---
$ws.Range("F11").Value = "<value>"
$ws.Range("C29").Value = "Small"
$ws.Range("D29").Value = "Right"

# --- New rows 31-40: TEXT ID, TYPOGRAPHY NAME, ALIGNMENT, DIRECTION, GB ---
$rows = @(
    @{ Row = 31; B = "off";            C = "LCD_Default"; D = "Left";   E = "LTR"; F = "off" },
    @{ Row = 32; B = "SingleUseId49";  C = "Small";        D = "Center"; E = "LTR"; F = "Reset" },
    @{ Row = 33; B = "SingleUseId50";  C = "Small";        D = "Left";   E = "LTR"; F = "Load" },
    @{ Row = 34; B = "SingleUseId51";  C = "Small";        D = "Left";   E = "LTR"; F = "A" },
    @{ Row = 35; B = "SingleUseId52";  C = "LCD_Default";  D = "Right";  E = "LTR"; F = "8.250<value>" },
    @{ Row = 36; B = "SingleUseId53";  C = "Tiny";         D = "Left";   E = "LTR"; F = "Current" },
    @{ Row = 37; B = "SingleUseId54";  C = "Small";        D = "Left";   E = "LTR"; F = "V" },
    @{ Row = 38; B = "SingleUseId55";  C = "LCD_Default";  D = "Right";  E = "LTR"; F = "25.1<value>" },
    @{ Row = 39; B = "SingleUseId56";  C = "Tiny";         D = "Left";   E = "LTR"; F = "Voltage" },
    @{ Row = 40; B = "SingleUseId57";  C = "Small";        D = "Center"; E = "LTR"; F = "Start" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
